{"js": "// Office.js (Word JavaScript API) edit script.\n// Body is the async (context) => { ... } function content.\n// Replaces the text of each of the 100 multiplication-problem cells\n// in the document's single table (20 rows x 5 columns, row-major order)\n// with the new values from the target revision.\n\nconst newValues = [\n  \"95\u00d754=\",\n  \"33\u00d753=\",\n  \"67\u00d749=\",\n  \"64\u00d797=\",\n  \"27\u00d794=\",\n  \"92\u00d733=\",\n  \"90\u00d731=\",\n  \"49\u00d793=\",\n  \"36\u00d786=\",\n  \"73\u00d764=\",\n  \"11\u00d711=\",\n  \"16\u00d714=\",\n  \"35\u00d769=\",\n  \"74\u00d725=\",\n  \"18\u00d748=\",\n  \"51\u00d773=\",\n  \"21\u00d729=\",\n  \"59\u00d751=\",\n  \"28\u00d710=\",\n  \"69\u00d798=\",\n  \"84\u00d738=\",\n  \"85\u00d721=\",\n  \"99\u00d790=\",\n  \"55\u00d773=\",\n  \"79\u00d747=\",\n  \"58\u00d740=\",\n  \"36\u00d738=\",\n  \"55\u00d712=\",\n  \"21\u00d757=\",\n  \"49\u00d752=\",\n  \"79\u00d783=\",\n  \"12\u00d790=\",\n  \"86\u00d770=\",\n  \"60\u00d713=\",\n  \"37\u00d725=\",\n  \"98\u00d761=\",\n  \"42\u00d768=\",\n  \"50\u00d712=\",\n  \"78\u00d774=\",\n  \"39\u00d771=\",\n  \"99\u00d715=\",\n  \"59\u00d755=\",\n  \"98\u00d779=\",\n  \"12\u00d722=\",\n  \"31\u00d722=\",\n  \"91\u00d773=\",\n  \"61\u00d711=\",\n  \"56\u00d720=\",\n  \"27\u00d723=\",\n  \"50\u00d782=\",\n  \"17\u00d718=\",\n  \"21\u00d713=\",\n  \"50\u00d795=\",\n  \"75\u00d724=\",\n  \"68\u00d7100=\",\n  \"50\u00d799=\",\n  \"64\u00d723=\",\n  \"62\u00d759=\",\n  \"52\u00d758=\",\n  \"85\u00d791=\",\n  \"84\u00d743=\",\n  \"80\u00d730=\",\n  \"66\u00d752=\",\n  \"48\u00d7100=\",\n  \"71\u00d725=\",\n  \"97\u00d730=\",\n  \"29\u00d784=\",\n  \"41\u00d720=\",\n  \"99\u00d753=\",\n  \"22\u00d714=\",\n  \"50\u00d794=\",\n  \"77\u00d785=\",\n  \"29\u00d757=\",\n  \"43\u00d737=\",\n  \"81\u00d758=\",\n  \"22\u00d718=\",\n  \"33\u00d765=\",\n  \"63\u00d748=\",\n  \"84\u00d725=\",\n  \"96\u00d764=\",\n  \"90\u00d746=\",\n  \"98\u00d788=\",\n  \"20\u00d758=\",\n  \"72\u00d790=\",\n  \"81\u00d768=\",\n  \"66\u00d739=\",\n  \"12\u00d7100=\",\n  \"74\u00d747=\",\n  \"83\u00d796=\",\n  \"24\u00d735=\",\n  \"76\u00d771=\",\n  \"24\u00d740=\",\n  \"72\u00d783=\",\n  \"60\u00d755=\",\n  \"61\u00d752=\",\n  \"95\u00d771=\",\n  \"80\u00d717=\",\n  \"28\u00d710=\",\n  \"11\u00d783=\",\n  \"21\u00d754=\"\n];\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nif (tables.items.length === 0) {\n  throw new Error(\"Expected a table in the document, found none.\");\n}\n\nconst table = tables.items[0];\ntable.load(\"values\");\nawait context.sync();\n\nconst rowCount = table.values.length;\nconst colCount = rowCount > 0 ? table.values[0].length : 0;\n\nif (rowCount * colCount !== newValues.length) {\n  throw new Error(\n    `Table shape (${rowCount}x${colCount}=${rowCount * colCount}) does not match expected cell count (${newValues.length}).`\n  );\n}\n\nlet i = 0;\nfor (let r = 0; r < rowCount; r++) {\n  for (let c = 0; c < colCount; c++) {\n    table.getCell(r, c).value = newValues[i];\n    i++;\n  }\n}\n\nawait context.sync();\n", "ps1": "# Word COM interop (PowerShell-style) edit script.\n# The document is open as $word.ActiveDocument. Replace the text of each of\n# the 100 multiplication-problem cells in the document's single table\n# (20 rows x 5 columns, row-major order) with the new values from the\n# target revision.\n\n$d = $word.ActiveDocument\n\n$newValues = @(\n    '95\u00d754=',\n    '33\u00d753=',\n    '67\u00d749=',\n    '64\u00d797=',\n    '27\u00d794=',\n    '92\u00d733=',\n    '90\u00d731=',\n    '49\u00d793=',\n    '36\u00d786=',\n    '73\u00d764=',\n    '11\u00d711=',\n    '16\u00d714=',\n    '35\u00d769=',\n    '74\u00d725=',\n    '18\u00d748=',\n    '51\u00d773=',\n    '21\u00d729=',\n    '59\u00d751=',\n    '28\u00d710=',\n    '69\u00d798=',\n    '84\u00d738=',\n    '85\u00d721=',\n    '99\u00d790=',\n    '55\u00d773=',\n    '79\u00d747=',\n    '58\u00d740=',\n    '36\u00d738=',\n    '55\u00d712=',\n    '21\u00d757=',\n    '49\u00d752=',\n    '79\u00d783=',\n    '12\u00d790=',\n    '86\u00d770=',\n    '60\u00d713=',\n    '37\u00d725=',\n    '98\u00d761=',\n    '42\u00d768=',\n    '50\u00d712=',\n    '78\u00d774=',\n    '39\u00d771=',\n    '99\u00d715=',\n    '59\u00d755=',\n    '98\u00d779=',\n    '12\u00d722=',\n    '31\u00d722=',\n    '91\u00d773=',\n    '61\u00d711=',\n    '56\u00d720=',\n    '27\u00d723=',\n    '50\u00d782=',\n    '17\u00d718=',\n    '21\u00d713=',\n    '50\u00d795=',\n    '75\u00d724=',\n    '68\u00d7100=',\n    '50\u00d799=',\n    '64\u00d723=',\n    '62\u00d759=',\n    '52\u00d758=',\n    '85\u00d791=',\n    '84\u00d743=',\n    '80\u00d730=',\n    '66\u00d752=',\n    '48\u00d7100=',\n    '71\u00d725=',\n    '97\u00d730=',\n    '29\u00d784=',\n    '41\u00d720=',\n    '99\u00d753=',\n    '22\u00d714=',\n    '50\u00d794=',\n    '77\u00d785=',\n    '29\u00d757=',\n    '43\u00d737=',\n    '81\u00d758=',\n    '22\u00d718=',\n    '33\u00d765=',\n    '63\u00d748=',\n    '84\u00d725=',\n    '96\u00d764=',\n    '90\u00d746=',\n    '98\u00d788=',\n    '20\u00d758=',\n    '72\u00d790=',\n    '81\u00d768=',\n    '66\u00d739=',\n    '12\u00d7100=',\n    '74\u00d747=',\n    '83\u00d796=',\n    '24\u00d735=',\n    '76\u00d771=',\n    '24\u00d740=',\n    '72\u00d783=',\n    '60\u00d755=',\n    '61\u00d752=',\n    '95\u00d771=',\n    '80\u00d717=',\n    '28\u00d710=',\n    '11\u00d783=',\n    '21\u00d754='\n)\n\n$table = $d.Tables.Item(1)\n$rowCount = $table.Rows.Count\n$colCount = $table.Columns.Count\n\nif ($rowCount * $colCount -ne $newValues.Count) {\n    throw \"Table shape ($rowCount x $colCount = $($rowCount * $colCount)) does not match expected cell count ($($newValues.Count)).\"\n}\n\n$i = 0\nfor ($r = 1; $r -le $rowCount; $r++) {\n    for ($c = 1; $c -le $colCount; $c++) {\n        $table.Cell($r, $c).Range.Text = $newValues[$i]\n        $i++\n    }\n}\n"}
